$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values such as "320.97" that look like plain
# numbers and would otherwise be silently reinterpreted by Excel as a
# number type. Only the specific D cells being updated need to be forced
# to text first, so each other cell keeps its original (default) style.
$priceCells = "D2","D3","D5","D6","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D22","D23","D24","D25","D26","D28","D29","D30","D33","D34","D36","D37","D38","D39","D40","D43","D45","D49","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "46.223.65"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.456.95"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "320.97"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "105.49"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "36.04"
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "18.37"
$ws.Range("E13").Value = "  -4.49%  "
$ws.Range("D14").Value = "7.09"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "2.843.59"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "2.468.04"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "46.124.53"
$ws.Range("E18").Value = "  +3.97%  "
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "70.97"
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "2.39"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").Value = "247.60"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").Value = "25.99"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "34.75"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +5.67%  "
$ws.Range("D33").Value = "19.75"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "0.0765"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.90"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("D40").Value = "125.25"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "20.97"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").Value = "1.973.47"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E48").Value = "  +11.88%  "
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("E50").Value = "  +9.89%  "
$ws.Range("D51").Value = "78.42"
$ws.Range("E51").Value = "  +5.15%  "
